# [IMP] New data for test environment
# Applies the withholding_tax.xlsx data update:
#  - rename a couple of wt codes (add trailing A/R markers)
#  - add a new "causale_pagamento_id" column (I) with external.A / external.R
#  - change wt_types (J) "other" -> "ritenuta"
#  - flip the Enasarco certification label from (A) to (R) and give it an
#    "external.R" causale_pagamento_id
#  - append a brand-new row 5 for the "1040-23%A" withholding tax

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: z0bug.wt_1040 -------------------------------------------------
$ws.Range("B2").Value = "1040-20%A"
$ws.Range("I2").Value = "external.A"
$ws.Range("J2").Value = "ritenuta"

# --- Row 3: z0bug.wt_1038 -------------------------------------------------
$ws.Range("B3").Value = "1040-23%R"
$ws.Range("I3").Value = "external.R"
$ws.Range("J3").Value = "ritenuta"

# --- Row 4: z0bug.wt_enasarco_1 -------------------------------------------
$ws.Range("C4").Value = "Enasarco 17% su 50% (R)"
$ws.Range("I4").Value = "external.R"

# --- Row 5 (new): z0bug.wt_1040-23A ---------------------------------------
$ws.Range("A5").Value = "z0bug.wt_1040-23A"
$ws.Range("B5").Value = "1040-23%A"
$ws.Range("C5").Value = "1040 – 23% su 100% (A)"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = $ws.Range("E4").Value2
$ws.Range("F5").Value = $ws.Range("F4").Value2
$ws.Range("G5").Value = $ws.Range("G4").Value2
$ws.Range("G5").Style = $ws.Range("G4").Style
$ws.Range("H5").Value = $ws.Range("H4").Value2
$ws.Range("I5").Value = "external.A"
$ws.Range("J5").Value = "ritenuta"
$ws.Range("K5").Value = 1

# Column B is now wider because of the longer codes
$ws.Columns("B").ColumnWidth = 12.7

# Match the author's final selection
$ws.Range("A2").Select() | Out-Null
